$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Sprint 6 rows (19-21) and the now "IN PROGRESS" status cell (C19) ---
# pick up the "Neutral" (with border) formatting that currently lives on A16
# (Sprint 5 group), before that cell's own style is changed below.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null

# Row 19 story ("I want to have place to store configurations...") moves
# from NOT STARTED to IN PROGRESS.
$ws.Range("C19").Value = "IN PROGRESS"

# --- Step 2: Sprint 5 rows (16-18) adopt the "Good" (no border) formatting ---
# already used elsewhere in column A (e.g. A2).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

# --- Step 3: row 18 ("I want to change app icon") becomes DONE ---
$ws.Range("C2").Copy() | Out-Null  # "Good" with border, matching other DONE cells
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "DONE"

$excel.CutCopyMode = 0

# --- Step 4: update active selection ---
$ws.Range("B16").Select()
